$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) contain numeric-looking / percent text that
# Excel would otherwise auto-convert; force them to Text format before writing.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.525.09'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.539.49'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.06%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '506.10'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.36'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -7.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.548.44'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.95%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.102'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.331'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.60%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.985.93'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.489.73'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.58'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.542.48'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.53'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -5.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '333.50'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -6.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.06'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.94'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.13'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.407'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.160'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.77%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.653.52'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.94%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0786'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -8.19%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.89'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -6.84%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.86'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.63%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '149.47'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.53'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.71%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.922'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.65%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.87'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.60%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.11'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.98%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.80'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.821'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -10.81%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -6.96%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '283.56'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.52'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0990'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.48%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.602'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.57%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0533'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.51%  '
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.29'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.60'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.52%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0227'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.63%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.52'
